$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date format (style) from an existing "Voltooid" date cell (D2, numFmtId 14)
# onto the two new date cells so they reuse the same cell style index instead of
# Excel minting a brand-new (duplicate) number-format/style entry.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D30:D31").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 30 - new "Verbetering uitleg en opmaak stored proc. en trigger" entry
$ws.Range("A30").Value = "Verbetering uitleg en opmaak stored proc. en trigger"
$ws.Range("B30").Value = "30 minuten"
$ws.Range("C30").Value = "1 uur"
$ws.Range("D30").Value = 41381
$ws.Range("E30").Value = 1
$ws.Range("F30").Value = "Steven V"
$ws.Range("G30").Value = "Solved"
$ws.Range("H30").Value = "DOC"

# Row 31 - new "Werken project (bugfixen)" entry
$ws.Range("A31").Value = "Werken project (bugfixen)"
$ws.Range("B31").Value = "6 uur"
$ws.Range("C31").Value = "6 uur"
$ws.Range("D31").Value = 41382
$ws.Range("E31").Value = 1
$ws.Range("F31").Value = "Wouter P, Robbie V, Steven V"
$ws.Range("G31").Value = "In Process"
$ws.Range("H31").Value = "APP"

# Column F widened because of the new, longer "Wouter P, Robbie V, Steven V" value.
$ws.Columns("F:F").ColumnWidth = 27

# Move / leave the selection where the author ended up after the edit.
$ws.Range("A32").Select() | Out-Null
